# "Sample Project" Rules workbook: the decision-table's last rule row label
# in column B (cell B11) is re-keyed from "R40" to "1". The new literal text
# is appended to the shared-string table and B11 is re-pointed at it.
#
# B11's existing cell style (the bottom-row border/fill of the decision
# table) must be left untouched, so "1" needs to land in B11 as TEXT without
# stamping a new style onto the cell. Assigning a numeric-looking string
# straight to .Value/.Value2 gets silently coerced to a number; forcing text
# via NumberFormat="@" (or a leading apostrophe) instead creates a brand-new
# style for the cell. To avoid both problems, the text value is produced in
# a scratch cell via a formula that evaluates to the text "1" (which keeps
# the scratch cell's own formatting untouched, i.e. General/no style bump),
# and only the resulting VALUE is brought across into B11 with a
# values-only paste, so B11 keeps its original style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")

$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""   # evaluates to the text string "1"

$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues: value only, keeps B11's own formatting

$scratch.Clear()
$excel.CutCopyMode = $false
